$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 77, shifting all
# existing rows (77-95) down by two (-> 79-97). This matches the
# diff: two brand-new price records are added, and everything that
# used to be row 77 onward moves down.
$ws.Range("A77:A78").EntireRow.Insert()

# New row 77: "Especial" quality record
$ws.Cells.Item(77,1).Value  = 4
$ws.Cells.Item(77,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(77,3).Value  = "Los Lagos"
$ws.Cells.Item(77,4).Value  = 44474
$ws.Cells.Item(77,5).Value  = 10
$ws.Cells.Item(77,6).Value  = "Fruta"
$ws.Cells.Item(77,7).Value  = 100101
$ws.Cells.Item(77,8).Value  = "Berries"
$ws.Cells.Item(77,9).Value  = 100112025
$ws.Cells.Item(77,10).Value = "Frutilla"
$ws.Cells.Item(77,11).Value = "Sin especificar"
$ws.Cells.Item(77,12).Value = "Especial"
$ws.Cells.Item(77,13).Value = 200
$ws.Cells.Item(77,14).Value = 15000
$ws.Cells.Item(77,15).Value = 15000
$ws.Cells.Item(77,16).Value = 15000
$ws.Cells.Item(77,17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(77,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(77,19).Value = 2143
$ws.Cells.Item(77,20).Value = 7

# New row 78: "Primera" quality record
$ws.Cells.Item(78,1).Value  = 4
$ws.Cells.Item(78,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(78,3).Value  = "Los Lagos"
$ws.Cells.Item(78,4).Value  = 44474
$ws.Cells.Item(78,5).Value  = 10
$ws.Cells.Item(78,6).Value  = "Fruta"
$ws.Cells.Item(78,7).Value  = 100101
$ws.Cells.Item(78,8).Value  = "Berries"
$ws.Cells.Item(78,9).Value  = 100112025
$ws.Cells.Item(78,10).Value = "Frutilla"
$ws.Cells.Item(78,11).Value = "Sin especificar"
$ws.Cells.Item(78,12).Value = "Primera"
$ws.Cells.Item(78,13).Value = 200
$ws.Cells.Item(78,14).Value = 13000
$ws.Cells.Item(78,15).Value = 13000
$ws.Cells.Item(78,16).Value = 13000
$ws.Cells.Item(78,17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(78,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(78,19).Value = 1857
$ws.Cells.Item(78,20).Value = 7
